$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.465.80"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.106.91"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.61"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5233"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4558"
$ws.Range("E8").Value = "  +5.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.27"
$ws.Range("E9").Value = "  +15.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08943"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.28"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "2.115.05"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.867"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.062"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.60"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001144"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.344"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "30.512.05"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.53"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").Value = "2.351.86"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.30"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.542"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.78"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.24"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("E31").Value = "  +2.47%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.669"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.381"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.944"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.49"
$ws.Range("E36").Value = "  +7.07%  "
$ws.Range("E37").Value = "  +5.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02587"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06852"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2305"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.73"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6892"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.252"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.324"
$ws.Range("E44").Value = "  +5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6385"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.253"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000346"
$ws.Range("E49").Value = "  +20.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3420"
$ws.Range("E50").Value = "  +25.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.53"
$ws.Range("E51").Value = "  +2.37%  "
